$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.341.27'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '2.228.81'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'244.43"
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = '  +2.15%  '
$ws.Range("D7").Value = "'73.80"
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = "'0.614"
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("D10").Value = "'42.49"
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("D11").Value = "'0.0977"
$ws.Range("E11").Value = '  +5.13%  '
$ws.Range("D12").Value = "'7.14"
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("E13").Value = '  +1.24%  '
$ws.Range("D14").Value = "'14.37"
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = "'0.853"
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = '2.224.72'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = "'0.0000117"
$ws.Range("E17").Value = '  +21.27%  '
$ws.Range("D18").Value = '42.123.54'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = "'6.16"
$ws.Range("E19").Value = '  +2.15%  '
$ws.Range("D20").Value = "'72.15"
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("D21").Value = "'10.02"
$ws.Range("E21").Value = '  +39.25%  '
$ws.Range("D22").Value = "'231.00"
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("E23").Value = '  -3.02%  '
$ws.Range("D24").Value = "'12.06"
$ws.Range("E24").Value = '  +9.93%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = "'3.63"
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("D27").Value = "'2.30"
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("E28").Value = '  +3.17%  '
$ws.Range("D29").Value = "'166.97"
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("D30").Value = "'21.06"
$ws.Range("E30").Value = '  +3.41%  '
$ws.Range("D31").Value = "'5.70"
$ws.Range("E31").Value = '  +19.14%  '
$ws.Range("D32").Value = "'0.0803"
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").Value = "'29.50"
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("D36").Value = "'4.41"
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("D37").Value = "'0.0304"
$ws.Range("E37").Value = '  +3.29%  '
$ws.Range("D38").Value = "'13.06"
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("D40").Value = "'5.62"
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").Value = "'62.43"
$ws.Range("E41").Value = '  +5.38%  '
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").Value = "'8.81"
$ws.Range("E43").Value = '  +2.68%  '
$ws.Range("D44").Value = "'105.32"
$ws.Range("E44").Value = '  -4.14%  '
$ws.Range("E45").Value = '  +3.34%  '
$ws.Range("D46").Value = "'0.994"
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = "'2.39"
$ws.Range("E47").Value = '  +8.74%  '
$ws.Range("E48").Value = '  +1.92%  '
$ws.Range("D49").Value = "'1.18"
$ws.Range("E49").Value = '  +2.57%  '
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").Value = "'4.05"
$ws.Range("E51").Value = '  +0.94%  '
